# ---------------------------------------------------------------------------
# "fix: ultima version del excel"
#
# 1) Hechos!B4/C4 and Hechos!B9/C9 relabel from the "PK" legend entry to a
#    new "Medida/Atributo" legend entry (same teal/cyan legend swatch style
#    already used by the "Media"/"Medida" cells in G/H).
#    Hechos!C9 specifically gets its OWN distinct new string "Media/Atributo".
#    Hechos!B8/C8 keep their text but pick up that same teal/cyan swatch
#    style too.
# 2) The merge-cell lists on ER and Dimensiones get re-sorted into reading
#    order (left-to-right). UnMerge()/Merge() resets the alignment/format
#    of the "filler" cells inside a merged range, so we snapshot + restore
#    those formats around the remerge to avoid collateral formatting
#    changes.
# 3) Dimensiones!A1 picks up the Arial header style already used by the
#    other top-level group headers (matches ER!G1 / ER!N1).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsER = $wb.Worksheets.Item("ER")
$wsHechos = $wb.Worksheets.Item("Hechos")
$wsDim = $wb.Worksheets.Item("Dimensiones")

# ---------------------------------------------------------------------------
# 1) Hechos: relabel legend cells and restyle them to the teal/cyan swatch
# ---------------------------------------------------------------------------

$wsHechos.Range("B4").Value = "Medida/Atributo"
$wsHechos.Range("C4").Value = "Medida/Atributo"

$wsHechos.Range("B9").Value = "Medida/Atributo"
$wsHechos.Range("C9").Value = "Media/Atributo"

# Copy the existing teal/cyan "Media" swatch formatting onto the cells that
# need it (B4, C4, B8, C8, B9, C9) without touching their values.
$wsHechos.Range("G4").Copy()
$wsHechos.Range("B4").PasteSpecial(-4122)
$wsHechos.Range("C4").PasteSpecial(-4122)
$wsHechos.Range("B8").PasteSpecial(-4122)
$wsHechos.Range("C8").PasteSpecial(-4122)
$wsHechos.Range("B9").PasteSpecial(-4122)
$wsHechos.Range("C9").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Re-sort merged ranges into left-to-right order on ER and Dimensiones.
#    Snapshot row-1 formatting to a scratch row first, since UnMerge()
#    clobbers the alignment/style of the cells it un-merges, then restore
#    it once the ranges are re-merged in the desired order.
# ---------------------------------------------------------------------------

# --- ER sheet : row 1, columns A:R ---
$wsER.Range("A1:R1").Copy()
$wsER.Range("A100").PasteSpecial(-4122)

$wsER.Range("A1:D1").UnMerge()
$wsER.Range("E1:F1").UnMerge()
$wsER.Range("G1:H1").UnMerge()
$wsER.Range("I1:M1").UnMerge()
$wsER.Range("N1:R1").UnMerge()

$wsER.Range("A1:D1").Merge()
$wsER.Range("E1:F1").Merge()
$wsER.Range("G1:H1").Merge()
$wsER.Range("I1:M1").Merge()
$wsER.Range("N1:R1").Merge()

$wsER.Range("A100:R100").Copy()
$wsER.Range("A1").PasteSpecial(-4122)
$wsER.Range("A100:R100").Delete()

# --- Dimensiones sheet : row 1, columns A:H ---
$wsDim.Range("A1:H1").Copy()
$wsDim.Range("A100").PasteSpecial(-4122)

$wsDim.Range("A1:B1").UnMerge()
$wsDim.Range("C1:D1").UnMerge()
$wsDim.Range("E1:F1").UnMerge()
$wsDim.Range("G1:H1").UnMerge()

$wsDim.Range("A1:B1").Merge()
$wsDim.Range("C1:D1").Merge()
$wsDim.Range("E1:F1").Merge()
$wsDim.Range("G1:H1").Merge()

$wsDim.Range("A100:H100").Copy()
$wsDim.Range("A1").PasteSpecial(-4122)
$wsDim.Range("A100:H100").Delete()

# ---------------------------------------------------------------------------
# 3) Dimensiones!A1: switch to the Arial header style (matches ER!G1/N1).
#    Done last so it isn't clobbered by the merge-order snapshot/restore
#    above.
# ---------------------------------------------------------------------------

$wsER.Range("G1").Copy()
$wsDim.Range("A1").PasteSpecial(-4122)
